# daily auto push: 2026-01-29 02:49 UTC
# Insert two new rows of data for 2026/01/29 right before the existing
# "2026/12/29" block (which currently starts at row 742), shifting the
# remaining rows down by 2 (old row 783 -> new row 785).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 742; this pushes former rows
# 742..783 down to 744..785 (and updates the sheet dimension).
$ws.Rows.Item(742).Insert()
$ws.Rows.Item(742).Insert()

# Column A holds the date as plain text (e.g. "2026/01/29"), matching the
# rest of the sheet. Force text format first so Excel doesn't silently
# convert the string into a date serial number.
$ws.Cells.Item(742, 1).NumberFormat = "@"
$ws.Cells.Item(743, 1).NumberFormat = "@"

# Fill in the two newly inserted rows with the new data points.
$ws.Cells.Item(742, 1).Value = "2026/01/29"
$ws.Cells.Item(742, 2).Value = "木"
$ws.Cells.Item(742, 3).Value = 7
$ws.Cells.Item(742, 4).Value = 201

$ws.Cells.Item(743, 1).Value = "2026/01/29"
$ws.Cells.Item(743, 2).Value = "木"
$ws.Cells.Item(743, 3).Value = 10
$ws.Cells.Item(743, 4).Value = 201
